$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Update column D (Pers.No./ID) values: 60000460 -> 60000468, etc. (rows with an explicit "D" id)
$dRows = @(6,7,8,9,10,16,17)
foreach ($r in $dRows) {
    $ws.Cells.Item($r, 4).Value = 60000468
}

$ws.Cells.Item(20, 4).Value = 60000469
$ws.Cells.Item(22, 4).Value = 60000470
$ws.Cells.Item(24, 4).Value = 60000471
$ws.Cells.Item(26, 4).Value = 60000472

# Update column E values: 312 -> 319, etc.
$eRows = @(11,12,13,14,15,18,19)
foreach ($r in $eRows) {
    $ws.Cells.Item($r, 5).Value = 319
}

$ws.Cells.Item(21, 5).Value = 320
$ws.Cells.Item(23, 5).Value = 321
$ws.Cells.Item(25, 5).Value = 322
$ws.Cells.Item(27, 5).Value = 323
